# Update the chemotherapy_regimen filter used throughout the "startup" sheet's
# Neo4j/Cypher query cells from "FEC (3 week cycles)" to
# "Other treatment given as part of a CTSU protocol" (perf-script filter update).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$oldValue = "FEC (3 week cycles)"
$newValue = "Other treatment given as part of a CTSU protocol"

# B2:B4 hold the per-tab (Cases/Samples/Files) queries, C2:C4 hold the shared
# StatQuery (counts) text - all four distinct query strings embed the same
# chemotherapy_regimen filter literal that needs updating.
$targetAddresses = @("B2", "C2", "B3", "C3", "B4", "C4")

foreach ($addr in $targetAddresses) {
    $cell = $ws.Range($addr)
    $text = $cell.Value2
    if ($text -ne $null -and $text.Contains($oldValue)) {
        $cell.Value2 = $text.Replace($oldValue, $newValue)
    }
}

# Reflect the cursor position left behind by the edit (last cell touched/selected).
$ws.Range("D4").Select()
